$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# Add new row 9 data
$ws.Range("A9").Value = "Is dit artikel nog op voorraad?"
$ws.Range("B9").Value = "mailmind.test@zohomail.eu"
$ws.Range("C9").Value = "Testmail #7: Is dit artikel nog op voorraad?"
$ws.Range("D9").Value = "Productinformatie"
$ws.Range("F9").Value = "2025-07-29 21:41:54"
$ws.Range("G9").Value = "Nee"
$ws.Range("H9").Value = "Ja"
$ws.Range("I9").Value = "Nee"
$ws.Range("J9").Value = "Nee"

# Extend conditional formatting ranges to include row 9
$colRanges = @("D2:D8", "G2:G8", "H2:H8", "I2:I8", "J2:J8")
$newRanges = @("D2:D9", "G2:G9", "H2:H9", "I2:I9", "J2:J9")
for ($idx = 0; $idx -lt $colRanges.Length; $idx++) {
    $fcs = $ws.Range($colRanges[$idx]).FormatConditions
    $target = $ws.Range($newRanges[$idx])
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($target)
    }
}

# Update Dashboard rows 4 and 5: swap category order, update counts
$dash.Range("A4").Value = "Productinformatie"
$dash.Range("B4").Value = 2
$dash.Range("A5").Value = "Bestelling / Levering"
$dash.Range("B5").Value = 1
